$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.424.66"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "'2.629.54"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D5").Value = "'596.09"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'152.94"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +3.43%  "
$ws.Range("D10").Value = "'5.85"
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("D11").Value = "'0.397"
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'28.13"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "'3.099.20"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "'64.312.25"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "  +10.13%  "
$ws.Range("D17").Value = "'2.605.90"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'12.30"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'4.78"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "'349.76"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "'7.11"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'67.72"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "'9.28"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'1.67"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "'8.32"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").Value = "'549.09"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'0.0₃0914"
$ws.Range("E31").Value = "  +7.67%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("D34").Value = "'5.53"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("D35").Value = "'6.24"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'0.422"
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("D37").Value = "'165.44"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "'20.12"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'168.99"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "'41.59"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("E44").Value = "  +5.05%  "
$ws.Range("D45").Value = "'23.33"
$ws.Range("E45").Value = "  +7.78%  "
$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  +12.61%  "
$ws.Range("D47").Value = "'0.0591"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "'19.36"
$ws.Range("E51").Value = "  -0.15%  "
